$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "BOM_Board1_PCB1_2025-11-27"

# Row 2 - 1uF capacitors (C1, C2)
$ws.Cells.Item(2,1).Value = "'1"
$ws.Cells.Item(2,2).Value = 2
$ws.Cells.Item(2,3).Value = "1uF"
$ws.Cells.Item(2,4).Value = "C1,C2"
$ws.Cells.Item(2,5).Value = "C0402"
$ws.Cells.Item(2,6).Value = "1uF"
$ws.Cells.Item(2,7).Value = "CL05A105KO5NNNC"
$ws.Cells.Item(2,8).Value = "SAMSUNG(三星)"
$ws.Cells.Item(2,9).Value = "C29266"
$ws.Cells.Item(2,10).Value = "LCSC"

# Row 3 - 100kOhm resistor (R1)
$ws.Cells.Item(3,1).Value = "'2"
$ws.Cells.Item(3,2).Value = 1
$ws.Cells.Item(3,3).Value = "100kΩ"
$ws.Cells.Item(3,4).Value = "R1"
$ws.Cells.Item(3,5).Value = "R0402"
$ws.Cells.Item(3,6).Value = "100kΩ"
$ws.Cells.Item(3,7).Value = "RC0402FR-07100KL"
$ws.Cells.Item(3,8).Value = "YAGEO(国巨)"
$ws.Cells.Item(3,9).Value = "C60491"
$ws.Cells.Item(3,10).Value = "LCSC"

# Row 4 - ESP32-S3-Zero (U1)
$ws.Cells.Item(4,1).Value = "'3"
$ws.Cells.Item(4,2).Value = 1
$ws.Cells.Item(4,3).Value = "ESP32-S3-Zero"
$ws.Cells.Item(4,4).Value = "U1"
$ws.Cells.Item(4,5).Value = "COMM-SMD_18P-P2.54-L23.5-W18.0-TL"
$ws.Cells.Item(4,6).Value = "'"
$ws.Cells.Item(4,7).Value = "ESP32-S3-Zero"
$ws.Cells.Item(4,8).Value = "'"
$ws.Cells.Item(4,9).Value = "C9900152785"
$ws.Cells.Item(4,10).Value = "LCSC"

# Row 5 - TPS22918DBVR (U2)
$ws.Cells.Item(5,1).Value = "'4"
$ws.Cells.Item(5,2).Value = 1
$ws.Cells.Item(5,3).Value = "TPS22918DBVR"
$ws.Cells.Item(5,4).Value = "U2"
$ws.Cells.Item(5,5).Value = "SOT-23-6_L2.9-W1.6-P0.95-LS2.8-BR"
$ws.Cells.Item(5,6).Value = "'"
$ws.Cells.Item(5,7).Value = "TPS22918DBVR"
$ws.Cells.Item(5,8).Value = "TI(德州仪器)"
$ws.Cells.Item(5,9).Value = "C131941"
$ws.Cells.Item(5,10).Value = "LCSC"

# Row 6 - ZX-SH1.0-5PWT (XBSMB)
$ws.Cells.Item(6,1).Value = "'5"
$ws.Cells.Item(6,2).Value = 1
$ws.Cells.Item(6,3).Value = "ZX-SH1.0-5PWT"
$ws.Cells.Item(6,4).Value = "XBSMB"
$ws.Cells.Item(6,5).Value = "CONN-SMD_5P-P1.00_MEGASTAR_ZX-SH1.0-5PWT"
$ws.Cells.Item(6,6).Value = "'"
$ws.Cells.Item(6,7).Value = "ZX-SH1.0-5PWT"
$ws.Cells.Item(6,8).Value = "Megastar(兆星)"
$ws.Cells.Item(6,9).Value = "C7430447"
$ws.Cells.Item(6,10).Value = "LCSC"

# Row 7 - trailing blank row (A7 only, empty text)
$ws.Cells.Item(7,1).Value = "'"
